$d = $word.ActiveDocument

# --- Edit 1: move the _GoBack bookmark from the "Parte V. Validaciones:" paragraph
#     down into the following (empty) paragraph. ---
$seg1Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="1B1B51A1" w14:textId="4D0EB3EC" w:rsidR="00061A32" w:rsidRPr="00D96381" w:rsidRDefault="00061A32" w:rsidP="00061A32"><w:pPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="FF6600"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="FF6600"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Parte V. Validaciones:</w:t></w:r></w:p><w:p w14:paraId="109930B8" w14:textId="77777777" w:rsidR="00697F60" w:rsidRDefault="00697F60" w:rsidP="00061A32"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="262626"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$p61 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Parte V. Validaciones:*") {
        $p61 = $p
        break
    }
}
if ($p61 -eq $null) {
    throw "Could not find 'Parte V. Validaciones:' paragraph"
}
$p62 = $p61.Next()

$rng1 = $d.Range($p61.Range.Start, $p62.Range.End)
[void]$rng1.InsertXML($seg1Xml)

# --- Edit 2: insert a new bullet paragraph right after the "... lenguaje php" item,
#     about validating forms from the Front-End (html/css). ---
$seg2Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="262626"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="262626"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Validar los Formularios de la aplicación desde el Front-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="262626"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>End</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="262626"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="262626"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>html</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="262626"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="262626"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$pPhp = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*desde Back-End, es decir, el lenguaje php*") {
        $pPhp = $p
        break
    }
}
if ($pPhp -eq $null) {
    throw "Could not find the '...lenguaje php' paragraph"
}

[void]$pPhp.Range.InsertParagraphAfter()
$newPara = $pPhp.Next()
[void]$newPara.Range.InsertXML($seg2Xml)

Write-Output "Edits applied."
